$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "CW3M C577" model label (row 8, col A) to "CW3M C579"
$ws.Range("A8").Value = "CW3M C579"

# Append a new row 9, duplicating the numeric results of row 8 under a new
# model label "CW3M C584" (same run family as row 8's "Demo_Baseline WRB 2010-18").
$ws.Range("A9").Value = "CW3M C584"
$ws.Range("B9").Value = "Demo_Baseline WRB 2010-18"
$ws.Range("C9").Value = "2010-18"

$ws.Range("D9").Value = 1181.5808646666667
$ws.Range("D9").NumberFormat = "0.00"

$ws.Range("E9").Value = 1612.6987305555554
$ws.Range("E9").NumberFormat = "0.00"

$ws.Range("F9").Value = 14.207868333333332
$ws.Range("F9").NumberFormat = "0.00"

$ws.Range("G9").Value = 52.671807666666659
$ws.Range("G9").NumberFormat = "0.00"

$ws.Range("H9").Value = 5.2579661111111111
$ws.Range("H9").NumberFormat = "0.00"

$ws.Range("I9").Value = 8.7714771111111105
$ws.Range("I9").NumberFormat = "0.00"

$ws.Range("J9").Value = 2.7833133333333331
$ws.Range("J9").NumberFormat = "0.00"

$ws.Range("K9").Value = 592.07319488888891
$ws.Range("K9").NumberFormat = "0.00"

$ws.Range("L9").Value = 43.587652666666663
$ws.Range("L9").NumberFormat = "0.00"

$ws.Range("M9").Value = 1035.3851454444443
$ws.Range("M9").NumberFormat = "0.00"

$ws.Range("N9").Value = 1200.5520154444446
$ws.Range("N9").NumberFormat = "0.00"

$ws.Range("O9").Value = 505160.02083333331
$ws.Range("O9").NumberFormat = "0"

$ws.Range("P9").Value = 286902.89236111112
$ws.Range("P9").NumberFormat = "0"

$ws.Range("Q9").Value = -0.80739277777777785
$ws.Range("Q9").NumberFormat = "0.00"

$ws.Range("R9").Value = -0.00028255555555555559
$ws.Range("R9").NumberFormat = "0.000000"

# Move the selection to A10, matching where the cursor lands after entering
# the new row of data.
$ws.Range("A10").Select()
